$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1)
$ws.Range("A1").Value = "Date"
$ws.Range("B1").Value = "Image"
$ws.Range("C1").Value = "Image text"
$ws.Range("D1").Value = "Heading"
$ws.Range("E1").Value = "Teaser"
$ws.Range("F1").Value = "Heading2"
$ws.Range("G1").Value = "Video"
$ws.Range("H1").Value = "Video text"
$ws.Range("I1").Value = "Signature"
$ws.Range("J1").Value = "Heading3"
$ws.Range("K1").Value = "Extra video 1"
$ws.Range("L1").Value = "Extra video text 1"
$ws.Range("M1").Value = "Extra video 2"
$ws.Range("N1").Value = "Extra video text 2"
$ws.Range("O1").Value = "Extra video 3"
$ws.Range("P1").Value = "Extra video text 3"
$ws.Range("Q1").Value = "Extra video 4"
$ws.Range("R1").Value = "Extra video text 4"

# Row 2
$ws.Range("A2").Value = [DateTime]::new(2001,12,8)
$ws.Range("B2").Value = "file:///C:/01%20naamisuvanto/naamisuvanto/template_test/images/kuulumiset/nokipannukahvit.jpg"
$ws.Range("C2").Value = "Nokkipannu kahvit"

# Row 3
$ws.Range("A3").Value = [DateTime]::new(2001,11,26)
$ws.Range("D3").Value = "Luonnonlohikannat kasvussa"

# Header style: bold font + fill color
$headerRange = $ws.Range("A1:R1")
$headerRange.Font.Bold = $true
$headerRange.Interior.ThemeColor = 5
$headerRange.Interior.TintAndShade = 0.79998168889431442

# Date format for A2:A3
$ws.Range("A2:A3").NumberFormat = "m/d/yyyy"

# Hyperlink on B2
$ws.Hyperlinks.Add($ws.Range("B2"), "file:///C:/01%20naamisuvanto/naamisuvanto/template_test/images/kuulumiset/nokipannukahvit.jpg")

# AutoFilter
$ws.Range("A1:R1").AutoFilter()

# Column widths
$ws.Columns.Item(1).ColumnWidth = 10.140625
$ws.Columns.Item(2).ColumnWidth = 94.28515625
$ws.Columns.Item(3).ColumnWidth = 17.85546875
$ws.Columns.Item(4).ColumnWidth = 26.85546875
$ws.Columns.Item(5).ColumnWidth = 9.140625
$ws.Columns.Item(6).ColumnWidth = 11.5703125
$ws.Columns.Item(7).ColumnWidth = 8.5703125
$ws.Columns.Item(8).ColumnWidth = 12.5703125
$ws.Columns.Item(9).ColumnWidth = 12.5703125
$ws.Columns.Item(10).ColumnWidth = 11.5703125
$ws.Columns.Item(11).ColumnWidth = 14.5703125
$ws.Columns.Item(12).ColumnWidth = 18.7109375
$ws.Columns.Item(13).ColumnWidth = 14.5703125
$ws.Columns.Item(14).ColumnWidth = 18.7109375
$ws.Columns.Item(15).ColumnWidth = 14.5703125
$ws.Columns.Item(16).ColumnWidth = 18.7109375
$ws.Columns.Item(17).ColumnWidth = 14.5703125
$ws.Columns.Item(18).ColumnWidth = 18.7109375

# Page setup
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Selection
$ws.Range("E10").Select()
